# Apply the "by Coach" worksheet edits:
#  - Toggle the Yes/No "Started" values (column C) for a set of rows.
#  - Reset the frozen-pane scroll position back to the top (A2) and
#    clear the explicit selection leftover from the previous save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Row -> new "Started" value (column C), toggled Yes<->No per the edit.
$updates = @{
    8  = "Yes"
    9  = "No"
    14 = "Yes"
    16 = "Yes"
    17 = "No"
    18 = "No"
    19 = "No"
    20 = "Yes"
    26 = "Yes"
    27 = "Yes"
    30 = "Yes"
    32 = "No"
    33 = "No"
    36 = "No"
    38 = "Yes"
    39 = "Yes"
    41 = "Yes"
    43 = "No"
    44 = "No"
    45 = "Yes"
    46 = "No"
    52 = "Yes"
    59 = "No"
    62 = "Yes"
    65 = "No"
    74 = "Yes"
    77 = "No"
    80 = "No"
    83 = "Yes"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Scroll the frozen pane back to the top and drop the saved C81 selection.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
